# Applies the "Se actualiza el código 11112020" commit:
#  - Tipo de Permiso value changes from "ANUENCIA DE FUNCIONAMIENTO" to "ALTA DE HACIENDA"
#  - A new "Factura" block (rows 45-51) is added, mirroring the existing "Recibo" block
#  - A new "Observaciones:" / "Estatus:" pair is added (rows 53-54)
#  - A new "Opciones tabla permisos:" block is added (rows 56-58) with A56:A58 merged
#    and centred alignment

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Factura block (mirrors the Recibo block above it) ----------------------
$ws.Range("A45").Value = "Factura"

$ws.Range("A46").Value = "Folio:"
$ws.Range("B46").Value = "F001"

$ws.Range("A47").Value = "Tipo Documento:"
$ws.Range("B47").Value = "ORIGINAL"

$ws.Range("A48").Value = "Importe:"
$ws.Range("B48").Value = "'1000"

$ws.Range("A49").Value = "Fecha:"
$ws.Range("B49").Value = "'21-09-2020"

$ws.Range("A50").Value = "Metodo Pago:"
$ws.Range("B50").Value = "ORIGINAL"

$ws.Range("A51").Value = "Ruta:"
$ws.Range("B51").Value = $ws.Range("B43").Value2

# --- Observaciones / Estatus --------------------------------------------------
$ws.Range("A53").Value = "Observaciones:"
$ws.Range("B53").Value = "Observaciones prueba"

$ws.Range("A54").Value = "Estatus:"
$ws.Range("B54").Value = "INGRESADO"

# --- Opciones tabla permisos: block (merged, centred label) -----------------
$ws.Range("A56").Value = "Opciones tabla permisos:"
$ws.Range("B56").Value = "detalle"

$ws.Range("B57").Value = "editar"

$ws.Range("B58").Value = "eliminar"

$ws.Range("A56:A58").Merge()
$ws.Range("A56:A58").HorizontalAlignment = -4108
$ws.Range("A56:A58").VerticalAlignment = -4108

# --- Tipo de Permiso: value update -----------------------------------------
$ws.Range("B20").Value = "ALTA DE HACIENDA"

# --- Selection / view state, best effort -------------------------------------
$ws.Range("B21").Select()
